# Fix typos in the testcase file:
#  - rename the "busbar" sheet to "bus"
#  - make the renamed "bus" sheet the active tab (it was "transformer" before)

$wb = $excel.ActiveWorkbook

$busSheet = $wb.Worksheets.Item("busbar")
$busSheet.Name = "bus"

# Activating this sheet both sets workbook bookViews/activeTab to it and
# marks its sheetView as the selected tab, implicitly clearing the
# previously-selected "transformer" tab.
$busSheet.Activate()
